$d = $word.ActiveDocument

# --- Change 1: merge "Jort Siemes" + " " runs into "Jort Siemes " ---
$d.Content.Find.Execute("Jort Siemes ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Jort Siemes ", 2) | Out-Null

Write-Output "done"
